$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.190.93'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '3.672.41'
$ws.Range("E3").Value = '  -3.56%  '
$ws.Range("D5").Value = '''596.49'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  -4.07%  '
$ws.Range("D7").Value = '3.670.31'
$ws.Range("E7").Value = '  -3.54%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("D11").Value = '''6.29'
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '''0.457'
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("D13").Value = '''37.81'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '4.286.28'
$ws.Range("E15").Value = '  -3.51%  '
$ws.Range("D16").Value = '3.675.29'
$ws.Range("E16").Value = '  -3.38%  '
$ws.Range("D17").Value = '68.115.26'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").Value = '''17.03'
$ws.Range("E20").Value = '  +5.85%  '
$ws.Range("D21").Value = '''491.00'
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '''9.07'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("D24").Value = '''84.33'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").Value = '''0.0000141'
$ws.Range("E25").Value = '  +2.14%  '
$ws.Range("D26").Value = '''2.28'
$ws.Range("E26").Value = '  -4.39%  '
$ws.Range("D27").Value = '''12.15'
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''7.82'
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''2.38'
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").Value = '''31.22'
$ws.Range("E33").Value = '  -4.92%  '
$ws.Range("D34").Value = '3.811.34'
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("D36").Value = '3.613.73'
$ws.Range("E36").Value = '  -3.55%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = '''0.993'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  -4.00%  '
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("D42").Value = '''428.96'
$ws.Range("E42").Value = '  -6.07%  '
$ws.Range("D43").Value = '''48.61'
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").Value = '''1.94'
$ws.Range("E44").Value = '  -2.93%  '
$ws.Range("D45").Value = '''2.80'
$ws.Range("E45").Value = '  -3.59%  '
$ws.Range("D46").Value = '''8.34'
$ws.Range("E46").Value = '  +0.75%  '
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = '''40.20'
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("D49").Value = '''141.13'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0348'
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.723.31'
$ws.Range("E51").Value = '  -3.79%  '
